$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "reportname" (rId1 / sheet1.xml): report identity now points
# at the "Order Detail" report instead of "Policy".
# ------------------------------------------------------------------
$sReport = $wb.Worksheets.Item("reportname")
$sReport.Range("A2").Value = "Order Detail"
$sReport.Range("A2").Select()

# ------------------------------------------------------------------
# Sheet "Sheets" (rId2 / sheet2.xml): only the remembered selection
# moves - no data changes.
# ------------------------------------------------------------------
$sSheets = $wb.Worksheets.Item("Sheets")
$sSheets.Range("C15").Select()

# ------------------------------------------------------------------
# Sheet "database" (rId3 / sheet3.xml): database type switches from
# Hive to Sql, and this tab becomes the active one.
# ------------------------------------------------------------------
$sDatabase = $wb.Worksheets.Item("database")
$sDatabase.Range("A2").Value = "Sql"

# ------------------------------------------------------------------
# Sheet "SQL" (rId4 / sheet4.xml): populate the connection-details
# table used for the new Sql Server data source.
# ------------------------------------------------------------------
$sSql = $wb.Worksheets.Item("SQL")

$sSql.Range("A1").Value = "Server"
$sSql.Range("B1").Value = "Database"
$sSql.Range("C1").Value = "Username"
$sSql.Range("D1").Value = "Password"
$sSql.Range("E1").Value = "Query"
$sSql.Range("A1:E1").Font.Bold = $true

$sSql.Range("A2").Value = "ACITGDPRWN01"
$sSql.Range("B2").Value = "gosales"
$sSql.Range("C2").Value = "sa"
$sSql.Range("D2").Value = "Pass1234$"

$query = @'
select [ORDER_DETAIL_CODE],[QUANTITY],[UNIT_SALE_PRICE]*[QUANTITY] as Revenue from  [gosales].[ORDER_DETAILS]
 where Order_Detail_Code in('1000001',
'1000002',
'1000013',
'1000014',
'1000015',
'1000016',
'1000017') order by [ORDER_DETAIL_CODE]
'@
$sSql.Range("E2").Value = $query
$sSql.Range("E2").WrapText = $true

$sSql.Columns.Item(3).ColumnWidth = 8.666666666666666
$sSql.Columns.Item(4).ColumnWidth = 8.166666666666666
$sSql.Columns.Item(5).ColumnWidth = 45

$sSql.Rows.Item(2).RowHeight = 158.4

# ------------------------------------------------------------------
# Make "database" the active tab, as it is now the sheet of interest
# for wiring up the new connection.
# ------------------------------------------------------------------
$sDatabase.Activate()
$sDatabase.Range("A2").Select()
